$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" summary text on sheet "Hoja1" ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = @"
Conversión del día 💰
✅ Dólar paralelo: 68

Binance
✅ 1000 Bs = 12.59 = 51057.53 pesos
✅ 51057.53 pesos = 12.54 = 970.45 Bs

Promedio competencia
✅ Tasa pesos: 20
✅ Tasa Bs: 20
✅ % Ganancia: 20%
"@

$ws1.Range("A1").Value = $newText

# --- Update the rate figures on sheet "tasas" ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 79.44
$ws2.Range("O10").Value = 4056.01
$ws2.Range("N12").Value = 4069.96
$ws2.Range("O12").Value = 77.358
